$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-30"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 05-30)"

# Update the June (row 6) total-column value with the new data point
$ws.Range("I6").Value = 110

# Update the Total row (row 14) total-column value
$ws.Range("I14").Value = 661
